$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37387
$ws.Range("D2").Value = 54083327
$ws.Range("C3").Value = 90300
$ws.Range("D3").Value = 132394724
$ws.Range("C4").Value = 30943
$ws.Range("D4").Value = 45828771
$ws.Range("C5").Value = 8636
$ws.Range("D5").Value = 12837235
$ws.Range("C6").Value = 1966
$ws.Range("D6").Value = 2921506
$ws.Range("C7").Value = 151
$ws.Range("D7").Value = 221593
$ws.Range("C11").Value = 40873
$ws.Range("D11").Value = 55488375
$ws.Range("C12").Value = 9565
$ws.Range("D12").Value = 13835326
$ws.Range("C13").Value = 25780
$ws.Range("D13").Value = 37809434
$ws.Range("C14").Value = 8279
$ws.Range("D14").Value = 12287763
$ws.Range("C15").Value = 2136
$ws.Range("D15").Value = 3176383
$ws.Range("C16").Value = 413
$ws.Range("D16").Value = 608623
$ws.Range("C19").Value = 10146
$ws.Range("D19").Value = 13447304
$ws.Range("C20").Value = 13291
$ws.Range("D20").Value = 19196352
$ws.Range("C21").Value = 31470
$ws.Range("D21").Value = 46189651
$ws.Range("C22").Value = 10178
$ws.Range("D22").Value = 15131897
$ws.Range("C23").Value = 2614
$ws.Range("D23").Value = 3888563
$ws.Range("C24").Value = 500
$ws.Range("D24").Value = 744092
$ws.Range("C26").Value = 11591
$ws.Range("D26").Value = 15492871
$ws.Range("C27").Value = 7577
$ws.Range("D27").Value = 10978038
$ws.Range("C28").Value = 22333
$ws.Range("D28").Value = 32780165
$ws.Range("C29").Value = 7763
$ws.Range("D29").Value = 11552802
$ws.Range("C30").Value = 1948
$ws.Range("D30").Value = 2906499
$ws.Range("C31").Value = 361
$ws.Range("D31").Value = 538915
$ws.Range("C33").Value = 8234
$ws.Range("D33").Value = 10882164
$ws.Range("C34").Value = 3188
$ws.Range("D34").Value = 4601190
$ws.Range("C35").Value = 7721
$ws.Range("D35").Value = 11276046
$ws.Range("C36").Value = 3149
$ws.Range("D36").Value = 4666461
$ws.Range("C37").Value = 819
$ws.Range("D37").Value = 1219823
$ws.Range("C38").Value = 160
$ws.Range("D38").Value = 238232
$ws.Range("C40").Value = 2424
$ws.Range("D40").Value = 3276990
$ws.Range("C41").Value = 17078
$ws.Range("D41").Value = 24699698
$ws.Range("C42").Value = 50717
$ws.Range("D42").Value = 74364674
$ws.Range("C43").Value = 18900
$ws.Range("D43").Value = 28075338
$ws.Range("C44").Value = 5574
$ws.Range("D44").Value = 8300978
$ws.Range("C45").Value = 1186
$ws.Range("D45").Value = 1769545
$ws.Range("C46").Value = 61
$ws.Range("D46").Value = 89568
$ws.Range("C49").Value = 16536
$ws.Range("D49").Value = 22035676
$ws.Range("C50").Value = 1974
$ws.Range("D50").Value = 2864316
$ws.Range("C51").Value = 6745
$ws.Range("D51").Value = 9918077
$ws.Range("C52").Value = 2317
$ws.Range("D52").Value = 3460418
$ws.Range("C53").Value = 746
$ws.Range("D53").Value = 1114305
$ws.Range("C54").Value = 181
$ws.Range("D54").Value = 268333
$ws.Range("C55").Value = 19
$ws.Range("D55").Value = 28500
$ws.Range("C56").Value = 6678
$ws.Range("D56").Value = 9201233
$ws.Range("C57").Value = 903
$ws.Range("D57").Value = 1325754
$ws.Range("C58").Value = 2261
$ws.Range("D58").Value = 3355113
$ws.Range("C59").Value = 908
$ws.Range("D59").Value = 1351501
$ws.Range("C60").Value = 310
$ws.Range("D60").Value = 464758
$ws.Range("C62").Value = 16
$ws.Range("D62").Value = 24000
$ws.Range("C63").Value = 1316
$ws.Range("D63").Value = 1857885
$ws.Range("C64").Value = 15224
$ws.Range("D64").Value = 21994545
$ws.Range("C65").Value = 44364
$ws.Range("D65").Value = 64930270
$ws.Range("C66").Value = 15614
$ws.Range("D66").Value = 23208214
$ws.Range("C67").Value = 4544
$ws.Range("D67").Value = 6768292
$ws.Range("C68").Value = 911
$ws.Range("D68").Value = 1354668
$ws.Range("C72").Value = 14982
$ws.Range("D72").Value = 19765325
$ws.Range("C73").Value = 50631
$ws.Range("D73").Value = 73688065
$ws.Range("C74").Value = 144220
$ws.Range("D74").Value = 212492303
$ws.Range("C75").Value = 62941
$ws.Range("D75").Value = 93793000
$ws.Range("C76").Value = 20087
$ws.Range("D76").Value = 30012611
$ws.Range("C77").Value = 4741
$ws.Range("D77").Value = 7083223
$ws.Range("C78").Value = 256
$ws.Range("D78").Value = 379170
$ws.Range("C84").Value = 50140
$ws.Range("D84").Value = 68281528
$ws.Range("C85").Value = 4535
$ws.Range("D85").Value = 6569691
$ws.Range("C86").Value = 11426
$ws.Range("D86").Value = 16787315
$ws.Range("C87").Value = 3841
$ws.Range("D87").Value = 5724406
$ws.Range("C88").Value = 1331
$ws.Range("D88").Value = 1988989
$ws.Range("C89").Value = 282
$ws.Range("D89").Value = 420512
$ws.Range("C92").Value = 5302
$ws.Range("D92").Value = 7130557
$ws.Range("C93").Value = 1563
$ws.Range("D93").Value = 2250035
$ws.Range("C94").Value = 5074
$ws.Range("D94").Value = 7475529
$ws.Range("C95").Value = 1924
$ws.Range("D95").Value = 2866446
$ws.Range("C96").Value = 681
$ws.Range("D96").Value = 1020460
$ws.Range("C97").Value = 176
$ws.Range("D97").Value = 263113
$ws.Range("C100").Value = 3471
$ws.Range("D100").Value = 4600143
$ws.Range("C101").Value = 587
$ws.Range("D101").Value = 874164
$ws.Range("C106").Value = 10676
$ws.Range("D106").Value = 15494495
$ws.Range("C107").Value = 29026
$ws.Range("D107").Value = 42650924
$ws.Range("C108").Value = 9725
$ws.Range("D108").Value = 14462213
$ws.Range("C109").Value = 2669
$ws.Range("D109").Value = 3979707
$ws.Range("C113").Value = 9712
$ws.Range("D113").Value = 12839606
$ws.Range("C114").Value = 30120
$ws.Range("D114").Value = 43442736
$ws.Range("C115").Value = 65659
$ws.Range("D115").Value = 96106627
$ws.Range("C116").Value = 21224
$ws.Range("D116").Value = 31541965
$ws.Range("C117").Value = 6013
$ws.Range("D117").Value = 8959326
$ws.Range("C118").Value = 1114
$ws.Range("D118").Value = 1664771
$ws.Range("C123").Value = 25593
$ws.Range("D123").Value = 34202435
$ws.Range("C124").Value = 35566
$ws.Range("D124").Value = 51340153
$ws.Range("C125").Value = 76059
$ws.Range("D125").Value = 111247443
$ws.Range("C126").Value = 23659
$ws.Range("D126").Value = 35117735
$ws.Range("C127").Value = 6332
$ws.Range("D127").Value = 9410504
$ws.Range("C128").Value = 1211
$ws.Range("D128").Value = 1800911
$ws.Range("C132").Value = 31376
$ws.Range("D132").Value = 41692948
$ws.Range("C133").Value = 13131
$ws.Range("D133").Value = 19010118
$ws.Range("C134").Value = 32106
$ws.Range("D134").Value = 47161708
$ws.Range("C135").Value = 11412
$ws.Range("D135").Value = 16956204
$ws.Range("C136").Value = 2938
$ws.Range("D136").Value = 4380504
$ws.Range("C137").Value = 486
$ws.Range("D137").Value = 722990
$ws.Range("C140").Value = 10741
$ws.Range("D140").Value = 14331981
$ws.Range("C141").Value = 34691
$ws.Range("D141").Value = 50108039
$ws.Range("C142").Value = 80536
$ws.Range("D142").Value = 118004163
$ws.Range("C143").Value = 24174
$ws.Range("D143").Value = 35921718
$ws.Range("C144").Value = 6335
$ws.Range("D144").Value = 9453208
$ws.Range("C145").Value = 1419
$ws.Range("D145").Value = 2110730
$ws.Range("C148").Value = 28896
$ws.Range("D148").Value = 39015661
